# Updates the "想去人数" (column F) counters on the "展览" (sheet1),
# "演出" (sheet2) and "全部类型" (sheet4) worksheets, matching the data
# refresh described in the commit "Update gh-pages to output generated
# at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsPerformance = $wb.Worksheets.Item("演出")
$wsAllTypes   = $wb.Worksheets.Item("全部类型")

# Row -> new value for the "展览" sheet (column F)
$exhibitionUpdates = @{
    4  = 27
    7  = 1137
    9  = 227
    10 = 322
    11 = 8002
    13 = 9496
    14 = 76
    17 = 466
    25 = 42
    29 = 1624
    30 = 27
    32 = 301
    35 = 337
    37 = 931
    38 = 6
    41 = 403
    43 = 267
    44 = 112
    45 = 273
}

foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row -> new value for the "演出" sheet (column F)
$performanceUpdates = @{
    15 = 48
    20 = 348
}

foreach ($row in $performanceUpdates.Keys) {
    $wsPerformance.Range("F$row").Value = $performanceUpdates[$row]
}

# Row -> new value for the "全部类型" sheet (column F)
$allTypesUpdates = @{
    8  = 27
    11 = 1137
    16 = 322
    17 = 8002
    18 = 9496
    19 = 76
    21 = 466
    25 = 42
    28 = 1624
    29 = 27
    31 = 301
    33 = 337
    36 = 931
    39 = 403
    40 = 48
    42 = 267
    43 = 112
    44 = 273
    48 = 348
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
